$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.09539999999999998
$ws.Range("E2").Value = 0.178
$ws.Range("I2").Value = [double]"2.725496231534034e-05"
$ws.Range("J2").Value = [double]"1.758716436197433e-05"
$ws.Range("K2").Value = 9.65
$ws.Range("L2").Value = 0.1446776611694153
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("W2").Value = 0.0760441292356186
$ws.Range("X2").Value = 0.04045678164711175
$ws.Range("Y2").Value = 0.03558734758850685
$ws.Range("Z2").Value = 0.4976092729159322
$ws.Range("AA2").Value = [double]"8.751536070815042e-06"
$ws.Range("AB2").Value = 0.03917622514009786
$ws.Range("AC2").Value = -0.03916747360402704
$ws.Range("AD2").Value = 7.26
$ws.Range("AE2").Value = 0.220910470067834
$ws.Range("AF2").Value = 7.480910470067834
$ws.Range("AG2").Value = 7.480910470067834
$ws.Range("AH2").Value = 0.05542198851686316
$ws.Range("AI2").Value = 0.04941779281706074
$ws.Range("AJ2").Value = 0.05542198851686316
$ws.Range("AK2").Value = 0.04941779281706074
$ws.Range("AN2").Value = 157.8260869565217
$ws.Range("AP2").Value = 162.6284884797355

# --- Row 3 ---
$ws.Range("D3").Value = 0.09539999999999998
$ws.Range("E3").Value = 0.178
$ws.Range("I3").Value = [double]"2.725496231534034e-05"
$ws.Range("J3").Value = [double]"1.758716436197433e-05"
$ws.Range("K3").Value = 9.65
$ws.Range("L3").Value = 0.1446776611694153
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("W3").Value = 0.0760441292356186
$ws.Range("X3").Value = 0.04045678164711175
$ws.Range("Y3").Value = 0.03558734758850685
$ws.Range("Z3").Value = 0.4976092729159322
$ws.Range("AA3").Value = [double]"8.751536070815042e-06"
$ws.Range("AB3").Value = 0.03917622514009786
$ws.Range("AC3").Value = -0.03916747360402704
$ws.Range("AD3").Value = 7.26
$ws.Range("AE3").Value = 0.220910470067834
$ws.Range("AF3").Value = 7.480910470067834
$ws.Range("AG3").Value = 7.480910470067834
$ws.Range("AH3").Value = 0.05542198851686316
$ws.Range("AI3").Value = 0.04941779281706074
$ws.Range("AJ3").Value = 0.05542198851686316
$ws.Range("AK3").Value = 0.04941779281706074
$ws.Range("AN3").Value = 157.8260869565217
$ws.Range("AP3").Value = 162.6284884797355
